$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pais")

# Update the "last updated" timestamp text in A1
$ws.Range("A1").Value = "Datos actualizados a 22 de Mayo de 2020 a las 18:05"

# Estados Unidos (row 4)
$ws.Range("B4").Value = 1626311
$ws.Range("C4").Value = 5409
$ws.Range("D4").Value = 383099
$ws.Range("E4").Value = 1146638
$ws.Range("G4").Value = 220
$ws.Range("H4").Value = 96574

# Alemania (row 11)
$ws.Range("B11").Value = 179381
$ws.Range("C11").Value = 360
$ws.Range("E11").Value = 12056
$ws.Range("G11").Value = 16
$ws.Range("H11").Value = 8325

# India (row 14)
$ws.Range("B14").Value = 124073
$ws.Range("C14").Value = 5847
$ws.Range("E14").Value = 69509
$ws.Range("G14").Value = 123
$ws.Range("H14").Value = 3707

# Polonia (row 35)
$ws.Range("B35").Value = 20619
$ws.Range("C35").Value = 476
$ws.Range("E35").Value = 10906
$ws.Range("G35").Value = 10
$ws.Range("H35").Value = 982

# Rumania (row 40)
$ws.Range("E40").Value = 5769
$ws.Range("G40").Value = 10
$ws.Range("H40").Value = 1166

# Argelia (row 56)
$ws.Range("B56").Value = 7918
$ws.Range("C56").Value = 190
$ws.Range("D56").Value = 4256
$ws.Range("E56").Value = 3080
$ws.Range("G56").Value = 7
$ws.Range("H56").Value = 582

# Luxemburgo (row 69)
$ws.Range("B69").Value = 3981
$ws.Range("C69").Value = 1
$ws.Range("D69").Value = 3748
$ws.Range("E69").Value = 124

# Irak (row 70)
$ws.Range("B70").Value = 3964
$ws.Range("C70").Value = 87
$ws.Range("D70").Value = 2532
$ws.Range("E70").Value = 1285
$ws.Range("G70").Value = 7
$ws.Range("H70").Value = 147

# Grecia (row 79)
$ws.Range("B79").Value = 2874
$ws.Range("C79").Value = 21
$ws.Range("E79").Value = 1331
$ws.Range("G79").Value = 1
$ws.Range("H79").Value = 169

# Row 112/113 swap: Republica de Chipre moves above Niger (updated case numbers)
$ws.Range("A112").Value = "Republica de Chipre"
$ws.Range("B112").Value = 927
$ws.Range("C112").Value = 4
$ws.Range("D112").Value = 561
$ws.Range("E112").Value = 349
$ws.Range("H112").Value = 17

$ws.Range("A113").Value = "Niger"
$ws.Range("B113").Value = 924
$ws.Range("D113").Value = 753
$ws.Range("E113").Value = 111
$ws.Range("H113").Value = 60

# Jordania (row 124)
$ws.Range("B124").Value = 700
$ws.Range("C124").Value = 16
$ws.Range("D124").Value = 461
$ws.Range("E124").Value = 230

# Etiopia (row 138)
$ws.Range("B138").Value = 433
$ws.Range("C138").Value = 34
$ws.Range("E138").Value = 300

# Libia (row 176)
$ws.Range("B176").Value = 72
$ws.Range("C176").Value = 1
$ws.Range("D176").Value = 38
$ws.Range("E176").Value = 31
